# feat: add 2022-Q1 data
#
# 1. A brand new sheet "2022-Q1" is inserted right after "2021-Q3" and
#    before "总计". It holds the per-fund holdings snapshot for the new
#    quarter (mirrors the layout of "2021-Q3").
# 2. The "总计" (totals) sheet gets a new row on top summarizing the
#    2022-Q1 quarter (1 holding, 0.01 亿元), pushing the existing
#    2021-Q3 totals row down.

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item(1)        # "2021-Q3"
$totalName = $wb.Worksheets.Item(2).Name   # "总计"

# --- Step 1: create the "2022-Q1" sheet -----------------------------------
# Duplicate the "总计" sheet (so the header/first-column cell styling,
# which uses the same style set as "总计", comes along for free) and place
# the copy right after "2021-Q3". (Copy shifts sheet positions, so the
# original "总计" sheet is re-fetched by name afterwards instead of by a
# now-stale index/reference.)
$wb.Worksheets.Item($totalName).Copy([System.Reflection.Missing]::Value, $q3)
$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"
$total = $wb.Worksheets.Item($totalName)

# Clear out the copied "总计" data (values only, keep the header/column
# styling that came along with the copy) before writing the new table.
$q1.Range("A1:D2").ClearContents()

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Re-apply the bold/bordered header style (style already lives on B1:D1
# from the copied sheet) across the full new header span.
$q1.Range("B1:D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Data row. Values such as "004098"/"0.34" must stay text (leading zeros /
# exact decimal strings), so they're entered with a leading apostrophe to
# force text entry instead of being auto-converted to numbers; the style
# is then reset to "Normal" so the apostrophe hint doesn't leave behind a
# quote-prefixed cell style.
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'004098"
$q1.Range("C2").Value = "前海开源港股通股息率50强股票"
$q1.Range("D2").Value = "'0.34"
$q1.Range("E2").Value = "'88.92"
$q1.Range("F2").Value = "'2.15"
$q1.Range("G2").Value = "'0.0073"
$q1.Range("H2").Value = 10
$q1.Range("B2").Style = "Normal"
$q1.Range("D2:G2").Style = "Normal"

# --- Step 2: insert the new totals row into "总计" --------------------------
# Push the existing data row (2021-Q3 totals) down to row 3, carrying its
# style along, then write the new 2022-Q1 summary row into row 2.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.13

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01
